# major accuracy check update
#
# 1) Sample number text "E7420" -> "E7420L" (shared string used by G2:G41)
# 2) Scroll the sheet view down a bit (topLeftCell A8 -> A12)
# 3) H2:H41 "accuracy check" cells: replace the =FALSE() formula with a
#    literal FALSE boolean value (same displayed/stored value, no formula)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename sample number E7420 -> E7420L across the whole column ---
$ws.Range("G2:G41").Value = "E7420L"

# --- 2) Update the remembered scroll position of the window ---
$win = $excel.ActiveWindow
$win.ScrollRow = 12
$win.ScrollColumn = 1

# --- 3) Replace the FALSE() formulas in H2:H41 with literal boolean FALSE ---
$ws.Range("H2:H41").Value = $false
